# Simulates a Word spell/grammar-check pass: wraps the flagged words
# ("prop", "Myeou", "Techau") with <w:proofErr w:type="spellStart/spellEnd"/>
# markers, and wraps the flagged phrase ("Computadora malo") with
# <w:proofErr w:type="gramStart/gramEnd"/> markers, splitting runs as needed.
#
# Word's object model has no direct API for inserting w:proofErr markers,
# so each affected paragraph's run(s) are rebuilt via Range.InsertXML with
# literal OOXML that reproduces the paragraph's original formatting plus
# the new proofErr markers / run split.

$d = $word.ActiveDocument
$ns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Replace-ParagraphXml($searchText, $xml, $matchWholeWord) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $matchWholeWord, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $searchText"
    }
    $null = $rng.InsertXML($xml)
}

# 1) "prop" -> spellStart/spellEnd around the whole run
Replace-ParagraphXml "prop" (
    "<w:p $ns w:rsidR='000F5B30' w:rsidRPr='00262F4D' w:rsidRDefault='000F5B30' w:rsidP='003E44BA'>" +
    "<w:pPr><w:rPr><w:sz w:val='12'/></w:rPr></w:pPr>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r w:rsidRPr='00262F4D'><w:rPr><w:sz w:val='12'/></w:rPr><w:t>prop</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "</w:p>"
), $true

# 2) "Myeou (vida completa)" -> split into "Myeou" (spellStart/spellEnd) + " (vida completa)"
Replace-ParagraphXml "Myeou (vida completa)" (
    "<w:p $ns w:rsidR='000F5B30' w:rsidRPr='00637FFE' w:rsidRDefault='000F5B30' w:rsidP='003E44BA'>" +
    "<w:pPr><w:ind w:left='2832' w:hanging='2832'/><w:rPr><w:sz w:val='16'/><w:u w:val='single'/></w:rPr></w:pPr>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:rPr><w:sz w:val='16'/></w:rPr><w:t>Myeou</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:rPr><w:sz w:val='16'/></w:rPr><w:t xml:space='preserve'> (vida completa)</w:t></w:r>" +
    "</w:p>"
), $false

# 3) "Myeou (media vida)" -> split into "Myeou" (spellStart/spellEnd) + " (media vida)"
Replace-ParagraphXml "Myeou (media vida)" (
    "<w:p $ns w:rsidR='000F5B30' w:rsidRPr='00262F4D' w:rsidRDefault='000F5B30' w:rsidP='003E44BA'>" +
    "<w:pPr><w:rPr><w:sz w:val='16'/></w:rPr></w:pPr>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:rPr><w:sz w:val='16'/></w:rPr><w:t>Myeou</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:rPr><w:sz w:val='16'/></w:rPr><w:t xml:space='preserve'> (media vida)</w:t></w:r>" +
    "</w:p>"
), $false

# 4) "Myeou (moribundo)" -> split into "Myeou" (spellStart/spellEnd) + " (moribundo)"
Replace-ParagraphXml "Myeou (moribundo)" (
    "<w:p $ns w:rsidR='000F5B30' w:rsidRPr='00262F4D' w:rsidRDefault='000F5B30' w:rsidP='003E44BA'>" +
    "<w:pPr><w:rPr><w:sz w:val='16'/></w:rPr></w:pPr>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:rPr><w:sz w:val='16'/></w:rPr><w:t>Myeou</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:rPr><w:sz w:val='16'/></w:rPr><w:t xml:space='preserve'> (moribundo)</w:t></w:r>" +
    "</w:p>"
), $false

# 5) standalone "Techau" -> spellStart/spellEnd around the whole run
Replace-ParagraphXml "Techau" (
    "<w:p $ns w:rsidR='000F5B30' w:rsidRDefault='000F5B30' w:rsidP='003E44BA'>" +
    "<w:pPr><w:rPr><w:sz w:val='16'/></w:rPr></w:pPr>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:rPr><w:sz w:val='16'/></w:rPr><w:t>Techau</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "</w:p>"
), $true

# 6) "Techau (" -> split into "Techau" (spellStart/spellEnd) + " (", keeping the
#    following "moribunda" and ")" runs of that paragraph untouched.
Replace-ParagraphXml "Techau (" (
    "<w:p $ns w:rsidR='000F5B30' w:rsidRDefault='000F5B30' w:rsidP='003E44BA'>" +
    "<w:pPr><w:rPr><w:sz w:val='16'/></w:rPr></w:pPr>" +
    "<w:proofErr w:type='spellStart'/>" +
    "<w:r><w:rPr><w:sz w:val='16'/></w:rPr><w:t>Techau</w:t></w:r>" +
    "<w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:rPr><w:sz w:val='16'/></w:rPr><w:t xml:space='preserve'> (</w:t></w:r>" +
    "<w:r w:rsidR='000D6589'><w:rPr><w:sz w:val='16'/></w:rPr><w:t>moribunda</w:t></w:r>" +
    "<w:r><w:rPr><w:sz w:val='16'/></w:rPr><w:t>)</w:t></w:r>" +
    "</w:p>"
), $false

# 7) "Computadora malo " -> split into "Computadora malo" (gramStart/gramEnd) + trailing space
Replace-ParagraphXml "Computadora malo " (
    "<w:p $ns w:rsidR='000F5B30' w:rsidRDefault='00767D01' w:rsidP='003E44BA'>" +
    "<w:pPr><w:rPr><w:sz w:val='16'/></w:rPr></w:pPr>" +
    "<w:proofErr w:type='gramStart'/>" +
    "<w:r><w:rPr><w:sz w:val='16'/></w:rPr><w:t>Computadora malo</w:t></w:r>" +
    "<w:proofErr w:type='gramEnd'/>" +
    "<w:r><w:rPr><w:sz w:val='16'/></w:rPr><w:t xml:space='preserve'> </w:t></w:r>" +
    "</w:p>"
), $false

Write-Output "proofErr edits applied"
